# Actualización semana 4 WEEK, Schedule, Task Personal
#
# Fills in the "Actual" figures for week 4 (row 17) of the schedule table
# on sheet "Hoja1": Team Hours (G17), Cumulative Hours (I17, formula),
# Week Earned Value (J17) and Cumulative Earned Value (L17, formula).
# Also moves the active selection to F17 (matching where the user was
# working when the week-4 Actual data was entered).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G17").Value = 5.5
$ws.Range("I17").Formula = "=G17+I16"
$ws.Range("J17").Value = 9.9
$ws.Range("L17").Formula = "=J17+L16"

[void]$ws.Range("F17").Select()
